{"js": "// Update the two-digit \u00f7 one-digit division worksheet: the table has 20\n// rows (5 populated with problems at rows 0, 4, 8, 12, 16, and 3 blank\n// rows after each for student work), 5 columns. Each populated cell's\n// text (\"NN\u00f7N=\") is replaced by a newly generated problem, per the\n// commit's regenerated-output diff. Cells are targeted by their\n// (row, column) position (not a blind text search) because several\n// \"before\"/\"after\" strings collide across cells (e.g. one cell's new\n// value equals another cell's old value), which would make ordinary\n// document-wide search-and-replace ambiguous/order-dependent.\nconst replacements = [\n  { row: 0, col: 0, before: \"73\u00f75=\", after: \"63\u00f76=\" },\n  { row: 0, col: 1, before: \"30\u00f75=\", after: \"82\u00f72=\" },\n  { row: 0, col: 2, before: \"14\u00f74=\", after: \"12\u00f73=\" },\n  { row: 0, col: 3, before: \"33\u00f76=\", after: \"61\u00f76=\" },\n  { row: 0, col: 4, before: \"59\u00f72=\", after: \"81\u00f75=\" },\n  { row: 4, col: 0, before: \"69\u00f75=\", after: \"54\u00f79=\" },\n  { row: 4, col: 1, before: \"22\u00f73=\", after: \"69\u00f79=\" },\n  { row: 4, col: 2, before: \"71\u00f77=\", after: \"78\u00f73=\" },\n  { row: 4, col: 3, before: \"64\u00f79=\", after: \"74\u00f79=\" },\n  { row: 4, col: 4, before: \"85\u00f76=\", after: \"84\u00f79=\" },\n  { row: 8, col: 0, before: \"49\u00f75=\", after: \"99\u00f79=\" },\n  { row: 8, col: 1, before: \"78\u00f73=\", after: \"62\u00f74=\" },\n  { row: 8, col: 2, before: \"77\u00f74=\", after: \"12\u00f73=\" },\n  { row: 8, col: 3, before: \"21\u00f77=\", after: \"86\u00f78=\" },\n  { row: 8, col: 4, before: \"26\u00f73=\", after: \"54\u00f78=\" },\n  { row: 12, col: 0, before: \"33\u00f72=\", after: \"12\u00f78=\" },\n  { row: 12, col: 1, before: \"44\u00f78=\", after: \"76\u00f79=\" },\n  { row: 12, col: 2, before: \"46\u00f76=\", after: \"12\u00f76=\" },\n  { row: 12, col: 3, before: \"35\u00f78=\", after: \"53\u00f75=\" },\n  { row: 12, col: 4, before: \"73\u00f78=\", after: \"14\u00f78=\" },\n  { row: 16, col: 0, before: \"31\u00f79=\", after: \"87\u00f72=\" },\n  { row: 16, col: 1, before: \"71\u00f75=\", after: \"63\u00f75=\" },\n  { row: 16, col: 2, before: \"19\u00f78=\", after: \"98\u00f74=\" },\n  { row: 16, col: 3, before: \"66\u00f78=\", after: \"64\u00f75=\" },\n  { row: 16, col: 4, before: \"33\u00f78=\", after: \"78\u00f72=\" },\n];\n\nconst table = context.document.body.tables.getFirst();\n\n// Grab every target cell's first paragraph up front (one sync for all\n// loads), then verify + rewrite each one.\nconst paragraphs = replacements.map(({ row, col }) =>\n  table.getCell(row, col).body.paragraphs.getFirst()\n);\nparagraphs.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nparagraphs.forEach((p, i) => {\n  const { before, after } = replacements[i];\n  const actual = p.text.trim();\n  if (actual !== before) {\n    throw new Error(\n      `Cell (${replacements[i].row}, ${replacements[i].col}) expected \"${before}\" but found \"${actual}\"`\n    );\n  }\n  p.insertText(after, Word.InsertLocation.replace);\n});\nawait context.sync();\n\n", "ps1": "# Update the two-digit \u00f7 one-digit division worksheet: the table has 20\n# rows (5 populated with problems at rows 1, 5, 9, 13, 17 in Word's\n# 1-based indexing, each followed by 3 blank rows for student work) and\n# 5 columns. Each populated cell's text (\"NN\u00f7N=\") is replaced by a\n# newly generated problem, per the commit's regenerated-output diff.\n#\n# Cells are targeted by their (row, column) position -- not a blind\n# Find/Replace-all -- because several \"before\"/\"after\" strings collide\n# across cells (e.g. one cell's new value equals another cell's old\n# value: \"71\u00f77=\" -> \"78\u00f73=\" while a different cell already holds\n# \"78\u00f73=\" -> \"62\u00f74=\"). A document-wide replace-all would be\n# ambiguous/order-dependent in that case.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; Before = '73\u00f75='; After = '63\u00f76=' }\n    @{ Row = 1; Col = 2; Before = '30\u00f75='; After = '82\u00f72=' }\n    @{ Row = 1; Col = 3; Before = '14\u00f74='; After = '12\u00f73=' }\n    @{ Row = 1; Col = 4; Before = '33\u00f76='; After = '61\u00f76=' }\n    @{ Row = 1; Col = 5; Before = '59\u00f72='; After = '81\u00f75=' }\n    @{ Row = 5; Col = 1; Before = '69\u00f75='; After = '54\u00f79=' }\n    @{ Row = 5; Col = 2; Before = '22\u00f73='; After = '69\u00f79=' }\n    @{ Row = 5; Col = 3; Before = '71\u00f77='; After = '78\u00f73=' }\n    @{ Row = 5; Col = 4; Before = '64\u00f79='; After = '74\u00f79=' }\n    @{ Row = 5; Col = 5; Before = '85\u00f76='; After = '84\u00f79=' }\n    @{ Row = 9; Col = 1; Before = '49\u00f75='; After = '99\u00f79=' }\n    @{ Row = 9; Col = 2; Before = '78\u00f73='; After = '62\u00f74=' }\n    @{ Row = 9; Col = 3; Before = '77\u00f74='; After = '12\u00f73=' }\n    @{ Row = 9; Col = 4; Before = '21\u00f77='; After = '86\u00f78=' }\n    @{ Row = 9; Col = 5; Before = '26\u00f73='; After = '54\u00f78=' }\n    @{ Row = 13; Col = 1; Before = '33\u00f72='; After = '12\u00f78=' }\n    @{ Row = 13; Col = 2; Before = '44\u00f78='; After = '76\u00f79=' }\n    @{ Row = 13; Col = 3; Before = '46\u00f76='; After = '12\u00f76=' }\n    @{ Row = 13; Col = 4; Before = '35\u00f78='; After = '53\u00f75=' }\n    @{ Row = 13; Col = 5; Before = '73\u00f78='; After = '14\u00f78=' }\n    @{ Row = 17; Col = 1; Before = '31\u00f79='; After = '87\u00f72=' }\n    @{ Row = 17; Col = 2; Before = '71\u00f75='; After = '63\u00f75=' }\n    @{ Row = 17; Col = 3; Before = '19\u00f78='; After = '98\u00f74=' }\n    @{ Row = 17; Col = 4; Before = '66\u00f78='; After = '64\u00f75=' }\n    @{ Row = 17; Col = 5; Before = '33\u00f78='; After = '78\u00f72=' }\n)\n\nforeach ($rep in $replacements) {\n    $cell = $t.Cell($rep.Row, $rep.Col)\n    $range = $cell.Range\n    # A table-cell Range.Text includes the trailing cell-mark characters;\n    # trim them off before comparing against the expected current value.\n    $current = $range.Text.TrimEnd([char]7, [char]13, [char]10)\n    if ($current -ne $rep.Before) {\n        throw \"Cell ($($rep.Row), $($rep.Col)) expected '$($rep.Before)' but found '$current'\"\n    }\n    $range.Text = $rep.After\n}\n\n"}
